# Generate Report for Handback
# Update the timestamp values recorded in the handback-status workbook to
# reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-08-29 01:03:58"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsZhCn.Range("H2").Value = "2016-08-29 01:03:53"
$wsZhCn.Range("K2").Value = "2016-08-29 01:04:13"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsDeDe.Range("H2").Value = "2016-08-29 01:03:58"
$wsDeDe.Range("K2").Value = "2016-08-29 01:04:21"
